# Initial commit for kfserving abtesting (#639)
#
# The big white "canvas" rectangle behind the mirroring diagram is
# resized/repositioned (its title/date-stamp header area above it grows),
# and a handful of now-unused helper shapes (stray empty textbox, leftover
# captions, a duplicate "v2.0" badge, the ribbon graphic, the "Objectives"
# checklist callout, and the footer/slide-number placeholders that had been
# copy-pasted onto this slide) are removed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Resize/reposition the big background rectangle (shape id 11).
$rect = $s.Shapes.Item("Rectangle 10")
$rect.Top    = 109.4115748031496
$rect.Height = 355.8066929133858

# Remove shapes that are no longer part of the slide.
$namesToDelete = @(
    "TextBox 9",                     # empty stray textbox (id 10)
    "TextBox 4",                     # "Istio Virtual Service" caption (id 5)
    "TextBox 95",                    # "winner" label (id 96)
    "Graphic 74",                    # Ribbon icon picture (id 75)
    "Rounded Rectangle 106",         # extra "v2.0" badge (id 107)
    "TextBox 124",                   # "Objectives" caption (id 125)
    "Graphic 125",                   # Checkbox Checked icon picture (id 126)
    "Footer Placeholder 2",          # footer placeholder (id 3)
    "Slide Number Placeholder 11"    # slide number placeholder (id 12)
)

foreach ($name in $namesToDelete) {
    $s.Shapes.Item($name).Delete()
}
